$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B8").Value = "13h 30m"
$ws.Range("B8").Select()
